$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New text for B1: comment.order_by (moves from D1 to B1, pushing out is_locked/is_enabled)
$ws.Range("B1").Value = "<%=comment.order_by%>"

# New text for C1: comment.rem (moves from E1 to C1)
$ws.Range("C1").Value = "<%=comment.rem%>"

# New text for D1: brand-new tenant_id_lbl validation column replacing the old order_by slot
$ws.Range("D1").Value = '<%=comment.tenant_id_lbl%><%selectList.tenant_id = data.findAllTenant.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.tenant_id.join(",") }"` })%>'

# E1 (old rem column) is removed entirely
$ws.Range("E1").ClearContents()

$wb.Save()
